$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 200005200
$ws.Range("I40").Value = 5997.5
$ws.Range("J40").Value = 333338000
$ws.Range("K40").Value = 5997.5
$ws.Range("L40").Value = 333338000
$ws.Range("M40").Value = -5822.5
$ws.Range("N40").Value = -333338350
$ws.Range("H62").Value = 2875.7058
$ws.Range("I62").Value = 2830.1428
$ws.Range("J62").Value = 3088.3333
$ws.Range("K62").Value = 2830.1428
$ws.Range("L62").Value = 3088.3333
$ws.Range("M62").Value = -2206.1428
$ws.Range("N62").Value = -4336.3333
$ws.Range("H65").Value = 2875.7058
$ws.Range("I65").Value = 2830.1428
$ws.Range("J65").Value = 3088.3333
$ws.Range("K65").Value = 14150.714
$ws.Range("L65").Value = 15441.6665
$ws.Range("M65").Value = -11030.714
$ws.Range("N65").Value = -21681.6665
$ws.Range("H76").Value = 8607
$ws.Range("I76").Value = 9229.833000000001
$ws.Range("J76").Value = 7361.3335
$ws.Range("K76").Value = 9229.833000000001
$ws.Range("L76").Value = 7361.3335
$ws.Range("M76").Value = -8914.833000000001
$ws.Range("N76").Value = -7991.3335
$ws.Range("H79").Value = 8607
$ws.Range("I79").Value = 9229.833000000001
$ws.Range("J79").Value = 7361.3335
$ws.Range("K79").Value = 9229.833000000001
$ws.Range("L79").Value = 7361.3335
$ws.Range("M79").Value = -8137.833000000001
$ws.Range("N79").Value = -9545.333500000001
$ws.Range("H86").Value = 2103.1738
$ws.Range("I86").Value = 1835.625
$ws.Range("K86").Value = 1835.625
$ws.Range("M86").Value = -712.625
$ws.Range("H89").Value = 2103.1738
$ws.Range("I89").Value = 1835.625
$ws.Range("K89").Value = 9178.125
$ws.Range("M89").Value = -3562.125
$ws.Range("H112").Value = 1635.9231
$ws.Range("I112").Value = 1095.909
$ws.Range("J112").Value = 2031.9333
$ws.Range("K112").Value = 3287.727
$ws.Range("L112").Value = 6095.7999
$ws.Range("M112").Value = -2179.727
$ws.Range("N112").Value = -8311.7999
$ws.Range("H127").Value = 2618.8235
$ws.Range("I127").Value = 2692.0625
$ws.Range("K127").Value = 8076.1875
$ws.Range("M127").Value = -3116.1875
$ws.Range("H132").Value = 6133.4443
$ws.Range("I132").Value = 6133.4443
$ws.Range("K132").Value = 18400.3329
$ws.Range("M132").Value = -15870.3329
$ws.Range("H138").Value = 2030.54
$ws.Range("I138").Value = 1139.7587
$ws.Range("J138").Value = 3260.6667
$ws.Range("K138").Value = 3419.2761
$ws.Range("L138").Value = 9782.000100000001
$ws.Range("M138").Value = 1720.7239
$ws.Range("N138").Value = -20062.0001
$ws.Range("H141").Value = 8715.352000000001
$ws.Range("I141").Value = 7971.4136
$ws.Range("J141").Value = 11412.125
$ws.Range("K141").Value = 23914.2408
$ws.Range("L141").Value = 34236.375
$ws.Range("M141").Value = -18734.2408
$ws.Range("N141").Value = -44596.375
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4462.522
$ws.Range("I32").Value = 3906.575
$ws.Range("K32").Value = 3906.575
$ws.Range("M32").Value = -3619.575
$ws.Range("H74").Value = 57840.78
$ws.Range("I74").Value = 42208.41
$ws.Range("J74").Value = 92232
$ws.Range("K74").Value = 42208.41
$ws.Range("L74").Value = 92232
$ws.Range("M74").Value = -41334.41
$ws.Range("N74").Value = -93980
$ws.Range("H77").Value = 57840.78
$ws.Range("I77").Value = 42208.41
$ws.Range("J77").Value = 92232
$ws.Range("K77").Value = 211042.05
$ws.Range("L77").Value = 461160
$ws.Range("M77").Value = -206674.05
$ws.Range("N77").Value = -469896
$ws.Range("H138").Value = 137500
$ws.Range("J138").Value = 137500
$ws.Range("L138").Value = 137500
$ws.Range("N138").Value = -147780
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 27532.111
$ws.Range("I102").Value = 4263
$ws.Range("J102").Value = 39166.668
$ws.Range("K102").Value = 4263
$ws.Range("L102").Value = 39166.668
$ws.Range("M102").Value = -1018
$ws.Range("N102").Value = -45656.668
$ws.Range("H134").Value = 3598.65
$ws.Range("I134").Value = 3797.2666
$ws.Range("K134").Value = 11391.7998
$ws.Range("M134").Value = -8856.799800000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8643
$ws.Range("I4").Value = 6667
$ws.Range("J4").Value = 10125
$ws.Range("K4").Value = 6667
$ws.Range("L4").Value = 10125
$ws.Range("M4").Value = -6555
$ws.Range("N4").Value = -10349
$ws.Range("H58").Value = 13272.5
$ws.Range("I58").Value = 7412.5
$ws.Range("J58").Value = 14444.5
$ws.Range("K58").Value = 7412.5
$ws.Range("L58").Value = 14444.5
$ws.Range("M58").Value = -7209.5
$ws.Range("N58").Value = -14850.5
$ws.Range("H134").Value = 8364.4
$ws.Range("I134").Value = 8182.6665
$ws.Range("K134").Value = 24547.9995
$ws.Range("M134").Value = -22012.9995
$ws.Range("H136").Value = 13272.5
$ws.Range("I136").Value = 7412.5
$ws.Range("J136").Value = 14444.5
$ws.Range("K136").Value = 22237.5
$ws.Range("L136").Value = 43333.5
$ws.Range("M136").Value = -19687.5
$ws.Range("N136").Value = -48433.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1780.7142
$ws.Range("J5").Value = 2555.7144
$ws.Range("L5").Value = 7667.1432
$ws.Range("N5").Value = -7891.1432
$ws.Range("H12").Value = 691.5
$ws.Range("J12").Value = 704.26666
$ws.Range("L12").Value = 2112.79998
$ws.Range("N12").Value = -2458.79998
$ws.Range("H14").Value = 2020
$ws.Range("I14").Value = 2020
$ws.Range("K14").Value = 6060
$ws.Range("M14").Value = -5887
$ws.Range("H64").Value = 8991.666999999999
$ws.Range("I64").Value = 4975
$ws.Range("J64").Value = 11000
$ws.Range("K64").Value = 14925
$ws.Range("L64").Value = 33000
$ws.Range("M64").Value = -14655
$ws.Range("N64").Value = -33540
$ws.Range("H67").Value = 8991.666999999999
$ws.Range("I67").Value = 4975
$ws.Range("J67").Value = 11000
$ws.Range("K67").Value = 14925
$ws.Range("L67").Value = 33000
$ws.Range("M67").Value = -13989
$ws.Range("N67").Value = -34872
$ws.Range("H114").Value = 1908.5
$ws.Range("I114").Value = 1784.5
$ws.Range("J114").Value = 2032.5
$ws.Range("K114").Value = 5353.5
$ws.Range("L114").Value = 6097.5
$ws.Range("M114").Value = -2099.5
$ws.Range("N114").Value = -12605.5
$ws.Range("H135").Value = 1780.7142
$ws.Range("J135").Value = 2555.7144
$ws.Range("L135").Value = 23001.4296
$ws.Range("N135").Value = -28071.4296
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 15797.875
$ws.Range("I39").Value = 16314
$ws.Range("J39").Value = 14249.5
$ws.Range("K39").Value = 16314
$ws.Range("L39").Value = 14249.5
$ws.Range("M39").Value = -15854
$ws.Range("N39").Value = -15169.5
$ws.Range("H132").Value = 2689.3333
$ws.Range("I132").Value = 2256
$ws.Range("J132").Value = 5072.6665
$ws.Range("K132").Value = 6768
$ws.Range("L132").Value = 15217.9995
$ws.Range("M132").Value = -4238
$ws.Range("N132").Value = -20277.9995
$ws.Range("H136").Value = 37727.69
$ws.Range("I136").Value = 2735.611
$ws.Range("K136").Value = 8206.832999999999
$ws.Range("M136").Value = -5656.832999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 26950.715
$ws.Range("J104").Value = 26950.715
$ws.Range("L104").Value = 26950.715
$ws.Range("N104").Value = -33938.715
$ws.Range("H126").Value = 43482130
$ws.Range("I126").Value = 58827256
$ws.Range("K126").Value = 176481768
$ws.Range("M126").Value = -176479298
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
